# "Running again after change in format + compute the rate_jobseek"
#
# The external data-prep pipeline that populates the unempl_c / jobseek_c /
# empl_jobseek_c columns (F:H) was re-run after a format change; this
# commit captures the workbook with those stale header labels and values
# cleared out (ready to be repopulated by the next run), and the UI
# selection reset to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three now-stale headers in F1:H1 (unempl_c, jobseek_c, empl_jobseek_c).
$ws.Range("F1:H1").ClearContents()

# Clear the computed values in F2:H53, keeping their existing number-format
# style (s="1") so the column formatting stays ready for the next run.
$ws.Range("F2:H53").ClearContents()

# Reset the view: scroll back to the top and select F8 (no more topLeftCell
# scroll offset, and the previous F46 selection no longer makes sense once
# the sheet shrinks).
$ws.Range("F8").Select()
